$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-15 (existing rows): columns A (index/number), B, C, E.
# Column D and F remain blank throughout (unchanged).
$existingRows = @(
    [ordered]@{ A=0; B="NSE:CONSOFINVT"; C="NSE:AARTECH"; E="NSE:BPCL" },
    [ordered]@{ A=1; B="NSE:GFLLIMITED"; C="NSE:ABSLAMC"; E="NSE:INDIAMART" },
    [ordered]@{ A=2; B="NSE:GOLDETF"; C="NSE:ADVENZYMES"; E="NSE:MRF" },
    [ordered]@{ A=3; B="NSE:GOLDSHARE"; C="NSE:ANGELONE"; E="NSE:PEL" },
    [ordered]@{ A=4; B="NSE:NBIFIN"; C="NSE:APOLSINHOT"; E="" },
    [ordered]@{ A=5; B=""; C="NSE:ARCHIDPLY"; E="" },
    [ordered]@{ A=6; B=""; C="NSE:ARVINDFASN"; E="" },
    [ordered]@{ A=7; B=""; C="NSE:ASIANENE"; E="" },
    [ordered]@{ A=8; B=""; C="NSE:AUSOMENT"; E="" },
    [ordered]@{ A=9; B=""; C="NSE:BIGBLOC"; E="" },
    [ordered]@{ A=10; B=""; C="NSE:CHEVIOT"; E="" },
    [ordered]@{ A=11; B=""; C="NSE:DBREALTY"; E="" },
    [ordered]@{ A=12; B=""; C="NSE:DEEPAKNTR"; E="" },
    [ordered]@{ A=13; B=""; C="NSE:DREDGECORP"; E="" }
)

$r = 2
foreach ($item in $existingRows) {
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 5).Value = $item.E
    $r = $r + 1
}

# New rows 16-36 appended below the existing data: only column A (index) and C (ticker) are populated;
# B, D, E, F stay blank. Column A needs the same bold/centered/bordered style as the rest of column A,
# so we copy formatting from an already-styled cell (A2) before writing each new value.
$newRows = @(
    [ordered]@{ A=14; C="NSE:EIHOTEL" },
    [ordered]@{ A=15; C="NSE:EIMCOELECO" },
    [ordered]@{ A=16; C="NSE:GOODLUCK" },
    [ordered]@{ A=17; C="NSE:GTECJAINX" },
    [ordered]@{ A=18; C="NSE:HTMEDIA" },
    [ordered]@{ A=19; C="NSE:IEX" },
    [ordered]@{ A=20; C="NSE:KFINTECH" },
    [ordered]@{ A=21; C="NSE:KHADIM" },
    [ordered]@{ A=22; C="NSE:MAGNUM" },
    [ordered]@{ A=23; C="NSE:MANALIPETC" },
    [ordered]@{ A=24; C="NSE:NOCIL" },
    [ordered]@{ A=25; C="NSE:NUVAMA" },
    [ordered]@{ A=26; C="NSE:NUVOCO" },
    [ordered]@{ A=27; C="NSE:PPL" },
    [ordered]@{ A=28; C="NSE:PRAJIND" },
    [ordered]@{ A=29; C="NSE:PRESTIGE" },
    [ordered]@{ A=30; C="NSE:RAILTEL" },
    [ordered]@{ A=31; C="NSE:RAMAPHO" },
    [ordered]@{ A=32; C="NSE:RAMRAT" },
    [ordered]@{ A=33; C="NSE:RAYMOND" },
    [ordered]@{ A=34; C="NSE:RVHL" }
)

$styleSource = $ws.Cells.Item(2, 1)
$r = 16
foreach ($item in $newRows) {
    $styleSource.Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 3).Value = $item.C
    $r = $r + 1
}

Write-Output "Updated sheet data through row 36"